$d = $word.ActiveDocument
$wdParagraph = 4
$wdCollapseStart = 1
$wdCollapseEnd = 0

function New-WordOpenXmlPackage([string]$bodyXml) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Replaces the contents of $paraRange (which must span a *whole* paragraph,
# end mark included) with the literal paragraph XML $newParaXml.
function Replace-ParagraphXml([object]$paraRange, [string]$newParaXml) {
    $xml = New-WordOpenXmlPackage($newParaXml)
    $paraRange.InsertXML($xml) | Out-Null
}

# Locates the Range covering the whole paragraph that contains $searchText.
# (Document.Paragraphs indexing becomes unreliable after InsertXML edits in
# this runtime, so we rely on Find + Range.Expand instead, which stay accurate.)
function Find-ParagraphRangeByText([string]$searchText) {
    $searchRange = $d.Content
    $found = $searchRange.Find.Execute($searchText)
    if (-not $found) {
        Write-Host "WARNING: text not found:" $searchText
        return $null
    }
    $searchRange.Collapse($wdCollapseStart) | Out-Null
    $searchRange.Expand($wdParagraph) | Out-Null
    return $searchRange
}

# --- 1) Table header row: wrap "toddlers", "teenagers", "adults", "elderly" with spell proofErr,
#        and split "toddlers " into "toddlers" + " " runs ---
$tbl = $d.Tables.Item(1)

$toddlersXml = '<w:p w:rsidR="00A54305" w:rsidRDefault="00A54305"><w:pPr><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>toddlers</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>'
Replace-ParagraphXml $tbl.Cell(1,1).Range $toddlersXml

$teenagersXml = '<w:p w:rsidR="00A54305" w:rsidRDefault="00A54305"><w:pPr><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>teenagers</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Replace-ParagraphXml $tbl.Cell(1,2).Range $teenagersXml

$adultsXml = '<w:p w:rsidR="00A54305" w:rsidRDefault="00A54305"><w:pPr><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>adults</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Replace-ParagraphXml $tbl.Cell(1,3).Range $adultsXml

$elderlyXml = '<w:p w:rsidR="00A54305" w:rsidRDefault="00A54305"><w:pPr><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>elderly</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Replace-ParagraphXml $tbl.Cell(1,4).Range $elderlyXml

# --- 2) Remove the bookmark around "miss" ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 3) Split the "…………………………when you are an adult or elderly to socialize with adolescents"
#        run into three runs: "…" + "down" + "………………………when you are an adult or elderly to socialize with adolescents" ---
$ellipsisParaRange = Find-ParagraphRangeByText("Sometimes it is advisable to lose your hair")
$ellipsisXml = '<w:p w:rsidR="00F629CD" w:rsidRPr="00DC12C1" w:rsidRDefault="00F629CD" w:rsidP="00F629CD"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00DC12C1"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Sometimes it </w:t></w:r><w:r w:rsidR="006A2962" w:rsidRPr="00DC12C1"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">is advisable to lose your hair </w:t></w:r><w:r w:rsidRPr="00DC12C1"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>…</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>down</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>………………………when you are an adult or elderly to socialize with adolescents</w:t></w:r></w:p>'
Replace-ParagraphXml $ellipsisParaRange $ellipsisXml

# --- 4) Change "loss" -> "desire" and add a new "_GoBack" bookmark right after it ---
$lossRange = $d.Content
$foundLoss = $lossRange.Find.Execute("loss…………………………..")
if ($foundLoss) {
    $lossRange.End = $lossRange.Start + 4
    $lossRange.Text = "desire"
    $bmStart = $lossRange.End
    $bmRange = $d.Range($bmStart, $bmStart)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
} else {
    Write-Host "WARNING: 'loss' run not found"
}

Write-Host "Done"
